# Reorders several MoSCoW-analysis bullet points:
#  - Moves "CRUD operations", the "Should Have" heading, "Cat Card - Edit Cats",
#    "Cat Food Types – Type", "Splash page – ads" and "Splash page – referrals"
#    paragraphs so they sit right after "Cat Feeding View – Basic" (before "DB Encryption").
#
# Word COM has no native "move paragraph" verb, so each paragraph is relocated
# by copying its FormattedText (preserving run/paragraph formatting such as the
# darkGreen highlight and the Heading1 style) to the new location and then
# deleting the original paragraph that got left behind.

$d = $word.ActiveDocument

function Find-ParaRange($doc, $text) {
    $rng = $doc.Content
    $ok = $rng.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw ("Find-ParaRange: text not found: " + $text)
    }
    $para = $rng.Paragraphs(1)
    return $para.Range
}

function Find-ParaRangeAfter($doc, $afterPos, $text) {
    $rng = $doc.Range($afterPos, $doc.Content.End)
    $ok = $rng.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw ("Find-ParaRangeAfter: text not found: " + $text)
    }
    $para = $rng.Paragraphs(1)
    return $para.Range
}

# --- Step 1: insert copies of the moved paragraphs right after "Cat Feeding
# View - Basic". Insert in reverse of the desired final order at a single
# fixed insertion point (immediately after the anchor) so each new paragraph
# pushes the previously-inserted ones down, ending up in the right order -
# this sidesteps the fact that setting .FormattedText does not advance the
# target Range's own End the way real Word would.

$anchorPara = Find-ParaRange $d "Cat Feeding View – Basic"
$insertAt = $anchorPara.End

$moveTextsReversed = @(
    "Splash page – referrals",
    "Splash page – ads",
    "Cat Food Types – Type",
    "Cat Card - Edit Cats",
    "Should Have",
    "CRUD operations"
)

foreach ($t in $moveTextsReversed) {
    $srcPara = Find-ParaRange $d $t
    $ft = $srcPara.FormattedText
    $ins = $d.Range($insertAt, $insertAt)
    $ins.FormattedText = $ft
}

# --- Step 2: delete the original paragraphs that were copied above. They
# now live further down, right after "DB Encryption" (which itself stays put
# - it simply ends up following the relocated block instead of preceding
# it). Scope every lookup to start after "DB Encryption" so the still-kept
# first copies (inserted in step 1) are never matched by accident.

$dbPara = Find-ParaRange $d "DB Encryption"
$afterDbPos = $dbPara.End

$dupTexts = @(
    "CRUD operations",
    "Should Have",
    "Cat Card - Edit Cats",
    "Cat Food Types – Type",
    "Splash page – ads",
    "Splash page – referrals"
)

foreach ($t in $dupTexts) {
    $dupPara = Find-ParaRangeAfter $d $afterDbPos $t
    $dupPara.Delete() | Out-Null
}

Write-Output "Reorder complete"
